$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 237.3
$ws.Range("I9").Value = 232.71428
$ws.Range("J9").Value = 248
$ws.Range("K9").Value = 232.71428
$ws.Range("L9").Value = 248
$ws.Range("M9").Value = -63.71428
$ws.Range("N9").Value = -586
$ws.Range("H19").Value = 3085
$ws.Range("I19").Value = 2145.889
$ws.Range("J19").Value = 4141.5
$ws.Range("K19").Value = 2145.889
$ws.Range("L19").Value = 4141.5
$ws.Range("M19").Value = -1970.889
$ws.Range("N19").Value = -4491.5
$ws.Range("H43").Value = 2061.3
$ws.Range("I43").Value = 2019.7142
$ws.Range("J43").Value = 2158.3333
$ws.Range("K43").Value = 2019.7142
$ws.Range("L43").Value = 2158.3333
$ws.Range("M43").Value = -1950.7142
$ws.Range("N43").Value = -2296.3333
$ws.Range("H98").Value = 6671.25
$ws.Range("I98").Value = 6671.25
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 6671.25
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -5173.25
$ws.Range("N98").ClearContents()
$ws.Range("H106").Value = 3143.1875
$ws.Range("I106").Value = 3391.6155
$ws.Range("K106").Value = 3391.6155
$ws.Range("M106").Value = -2760.6155
$ws.Range("H107").Value = 4753.027
$ws.Range("I107").Value = 4617.375
$ws.Range("J107").Value = 5621.2
$ws.Range("K107").Value = 4617.375
$ws.Range("L107").Value = 5621.2
$ws.Range("M107").Value = -2697.375
$ws.Range("N107").Value = -9461.200000000001
$ws.Range("H122").Value = 6671.25
$ws.Range("I122").Value = 6671.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 20013.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -17563.75
$ws.Range("N122").ClearContents()
$ws.Range("H127").Value = 1430.6923
$ws.Range("I127").Value = 1339.0435
$ws.Range("J127").Value = 2133.3333
$ws.Range("K127").Value = 4017.1305
$ws.Range("L127").Value = 6399.999899999999
$ws.Range("M127").Value = 942.8694999999998
$ws.Range("N127").Value = -16319.9999
$ws.Range("H138").Value = 2438.551
$ws.Range("J138").Value = 2973.4707
$ws.Range("L138").Value = 8920.4121
$ws.Range("N138").Value = -19200.4121
$ws.Range("H141").Value = 4149.769
$ws.Range("I141").Value = 4303.8184
$ws.Range("K141").Value = 12911.4552
$ws.Range("M141").Value = -7731.4552

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 37319.2
$ws.Range("I28").Value = 17238.4
$ws.Range("K28").Value = 17238.4
$ws.Range("M28").Value = -17046.4
$ws.Range("H74").Value = 7937072
$ws.Range("I74").Value = 8547597
$ws.Range("K74").Value = 8547597
$ws.Range("M74").Value = -8546723
$ws.Range("H77").Value = 7937072
$ws.Range("I77").Value = 8547597
$ws.Range("K77").Value = 42737985
$ws.Range("M77").Value = -42733617
$ws.Range("H99").Value = 37319.2
$ws.Range("I99").Value = 17238.4
$ws.Range("K99").Value = 17238.4
$ws.Range("M99").Value = -14243.4
$ws.Range("H110").Value = 3369.7932
$ws.Range("I110").Value = 3065.5
$ws.Range("K110").Value = 3065.5
$ws.Range("M110").Value = -1020.5
$ws.Range("H122").Value = 1861
$ws.Range("I122").Value = 1628.2941
$ws.Range("K122").Value = 4884.8823
$ws.Range("M122").Value = -2434.8823
$ws.Range("H132").Value = 4560.222
$ws.Range("I132").Value = 4877.657
$ws.Range("K132").Value = 14632.971
$ws.Range("M132").Value = -12102.971

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 541.26666
$ws.Range("J80").Value = 346.44446
$ws.Range("L80").Value = 346.44446
$ws.Range("N80").Value = -2342.44446
$ws.Range("H83").Value = 541.26666
$ws.Range("J83").Value = 346.44446
$ws.Range("L83").Value = 1732.2223
$ws.Range("N83").Value = -11716.2223
$ws.Range("H99").Value = 2919.8
$ws.Range("I99").Value = 2919.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2919.8
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1421.8
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 136368100
$ws.Range("I105").Value = 214287700
$ws.Range("K105").Value = 214287700
$ws.Range("M105").Value = -214285953
$ws.Range("H134").Value = 5192.4
$ws.Range("I134").Value = 4284.5264
$ws.Range("J134").Value = 8067.3335
$ws.Range("K134").Value = 12853.5792
$ws.Range("L134").Value = 24202.0005
$ws.Range("M134").Value = -10318.5792
$ws.Range("N134").Value = -29272.0005

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2059.913
$ws.Range("I58").Value = 2427.7144
$ws.Range("K58").Value = 2427.7144
$ws.Range("M58").Value = -2224.7144
$ws.Range("H99").Value = 3801.7273
$ws.Range("I99").Value = 3098.8
$ws.Range("J99").Value = 4387.5
$ws.Range("K99").Value = 3098.8
$ws.Range("L99").Value = 4387.5
$ws.Range("M99").Value = -1600.8
$ws.Range("N99").Value = -7383.5
$ws.Range("H126").Value = 3801.7273
$ws.Range("I126").Value = 3098.8
$ws.Range("J126").Value = 4387.5
$ws.Range("K126").Value = 9296.400000000001
$ws.Range("L126").Value = 13162.5
$ws.Range("M126").Value = -6826.400000000001
$ws.Range("N126").Value = -18102.5
$ws.Range("H132").Value = 2832.4644
$ws.Range("I132").Value = 2600.4783
$ws.Range("J132").Value = 3899.6
$ws.Range("K132").Value = 7801.4349
$ws.Range("L132").Value = 11698.8
$ws.Range("M132").Value = -5271.4349
$ws.Range("N132").Value = -16758.8
$ws.Range("H134").Value = 4008.1538
$ws.Range("I134").Value = 3516
$ws.Range("K134").Value = 10548
$ws.Range("M134").Value = -8013
$ws.Range("H136").Value = 2059.913
$ws.Range("I136").Value = 2427.7144
$ws.Range("K136").Value = 7283.1432
$ws.Range("M136").Value = -4733.1432

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 12155
$ws.Range("I3").Value = 9813.546
$ws.Range("J3").Value = 25033
$ws.Range("K3").Value = 29440.638
$ws.Range("L3").Value = 75099
$ws.Range("M3").Value = -29328.638
$ws.Range("N3").Value = -75323
$ws.Range("H4").Value = 139059660
$ws.Range("I4").Value = 166833360
$ws.Range("J4").Value = 111285960
$ws.Range("K4").Value = 500500080
$ws.Range("L4").Value = 333857880
$ws.Range("M4").Value = -500499968
$ws.Range("N4").Value = -333858104
$ws.Range("H68").Value = 3473600.2
$ws.Range("J68").Value = 5436152.5
$ws.Range("L68").Value = 16308457.5
$ws.Range("N68").Value = -16310079.5
$ws.Range("H71").Value = 3473600.2
$ws.Range("J71").Value = 5436152.5
$ws.Range("L71").Value = 48925372.5
$ws.Range("N71").Value = -48933484.5
$ws.Range("H124").Value = 24999
$ws.Range("I124").Value = 24999
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 74997
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("M124").Value = -70087
$ws.Range("H128").Value = 454931
$ws.Range("I128").Value = 454931
$ws.Range("K128").Value = 1364793
$ws.Range("M128").Value = -1359813

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8657.929
$ws.Range("I102").Value = 10465.375
$ws.Range("J102").Value = 6248
$ws.Range("K102").Value = 10465.375
$ws.Range("L102").Value = 6248
$ws.Range("M102").Value = -8843.375
$ws.Range("N102").Value = -9492
$ws.Range("H134").Value = 47721.168
$ws.Range("J134").Value = 47721.168
$ws.Range("L134").Value = 143163.504
$ws.Range("N134").Value = -148233.504

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1407.3334
$ws.Range("I16").Value = 1434.5385
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 1434.5385
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = -1264.5385
$ws.Range("N16").Value = -1040
$ws.Range("H40").Value = 6927.1816
$ws.Range("I40").Value = 6375.5
$ws.Range("J40").Value = 8398.333000000001
$ws.Range("K40").Value = 6375.5
$ws.Range("L40").Value = 8398.333000000001
$ws.Range("M40").Value = -6239.5
$ws.Range("N40").Value = -8670.333000000001
$ws.Range("H68").Value = 2680.353
$ws.Range("J68").Value = 4620
$ws.Range("L68").Value = 4620
$ws.Range("N68").Value = -6118
$ws.Range("H71").Value = 2680.353
$ws.Range("J71").Value = 4620
$ws.Range("L71").Value = 23100
$ws.Range("N71").Value = -30588
$ws.Range("H132").Value = 1980609.9
$ws.Range("I132").Value = 2168401.5
$ws.Range("K132").Value = 6505204.5
$ws.Range("M132").Value = -6502674.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 31666
$ws.Range("J41").Value = 31666
$ws.Range("L41").Value = 31666
$ws.Range("N41").Value = -32446
$ws.Range("H122").Value = 5311
$ws.Range("I122").Value = 4849.5386
$ws.Range("J122").Value = 7310.6665
$ws.Range("K122").Value = 14548.6158
$ws.Range("L122").Value = 21931.9995
$ws.Range("M122").Value = -12098.6158
$ws.Range("N122").Value = -26831.9995
$ws.Range("H126").Value = 113758410
$ws.Range("I126").Value = 113758410
$ws.Range("K126").Value = 341275230
$ws.Range("M126").Value = -341272760
$ws.Range("H135").Value = 70574.266
$ws.Range("J135").Value = 70574.266
$ws.Range("L135").Value = 70574.266
$ws.Range("N135").Value = -80714.266
